# This workbook is a "stacked weekly" price sheet: each week a brand-new
# observation is inserted as the new row 8 (right after the fixed top rows
# 1-7), pushing every existing data row down by one and growing the used
# range by a row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 - this shifts rows 8..51 down to
# 9..52 (and all their formatting/values) exactly like using Excel's
# "Insert Sheet Rows" command, and extends the sheet dimension to R52.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with this week's observation.
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 45163
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100112026
$ws.Cells.Item(8, 7).Value = "Haba"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 600
$ws.Cells.Item(8, 11).Value = 10000
$ws.Cells.Item(8, 12).Value = 11000
$ws.Cells.Item(8, 13).Value = 10500
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(8, 16).Value = 420
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
